$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2175732217573222
$ws.Range("C2").Value = 0.4560669456066946
$ws.Range("J2").Value = 0.01255230125523013
$ws.Range("P2").Value = 0.200836820083682
$ws.Range("S2").Value = 0.1129707112970711
$ws.Range("B3").Value = 0.01834862385321101
$ws.Range("C3").Value = 0.02752293577981652
$ws.Range("J3").Value = 0.02752293577981652
$ws.Range("O3").Value = 0.009174311926605505
$ws.Range("P3").Value = 0.7981651376146789
$ws.Range("S3").Value = 0.1192660550458716
$ws.Range("J4").Value = 0.05714285714285714
$ws.Range("O4").Value = 0.02857142857142857
$ws.Range("P4").Value = 0.6857142857142857
$ws.Range("S4").Value = 0.2285714285714286
$ws.Range("B6").Value = 0.07386363636363637
$ws.Range("D6").Value = 0.005681818181818182
$ws.Range("F6").Value = 0.06818181818181818
$ws.Range("J6").Value = 0.2386363636363636
$ws.Range("O6").Value = 0.01136363636363636
$ws.Range("Q6").Value = 0.1534090909090909
$ws.Range("R6").Value = 0.05681818181818182
$ws.Range("S6").Value = 0.3920454545454545
$ws.Range("B7").Value = 0.1148325358851675
$ws.Range("D7").Value = 0.009569377990430622
$ws.Range("E7").Value = 0.009569377990430622
$ws.Range("F7").Value = 0.02870813397129187
$ws.Range("J7").Value = 0.1339712918660287
$ws.Range("Q7").Value = 0.2200956937799043
$ws.Range("R7").Value = 0.0430622009569378
$ws.Range("S7").Value = 0.4401913875598086
$ws.Range("B8").Value = 0.07246376811594203
$ws.Range("D8").Value = 0.02484472049689441
$ws.Range("F8").Value = 0.06211180124223602
$ws.Range("J8").Value = 0.113871635610766
$ws.Range("O8").Value = 0.010351966873706
$ws.Range("Q8").Value = 0.1801242236024845
$ws.Range("R8").Value = 0.08281573498964803
$ws.Range("S8").Value = 0.453416149068323
$ws.Range("B9").Value = 0.07199999999999999
$ws.Range("F9").Value = 0.07199999999999999
$ws.Range("J9").Value = 0.112
$ws.Range("O9").Value = 0.016
$ws.Range("Q9").Value = 0.2
$ws.Range("R9").Value = 0.064
$ws.Range("S9").Value = 0.464
$ws.Range("B10").Value = 0.0926605504587156
$ws.Range("D10").Value = 0.01926605504587156
$ws.Range("F10").Value = 0.06422018348623854
$ws.Range("J10").Value = 0.1119266055045872
$ws.Range("O10").Value = 0.009174311926605505
$ws.Range("Q10").Value = 0.2330275229357798
$ws.Range("R10").Value = 0.06972477064220184
$ws.Range("S10").Value = 0.4
$ws.Range("G11").Value = 0.134020618556701
$ws.Range("J11").Value = 0.07216494845360824
$ws.Range("K11").Value = 0.1855670103092784
$ws.Range("L11").Value = 0.5876288659793815
$ws.Range("S11").Value = 0.02061855670103093
$ws.Range("G12").Value = 0.7954545454545454
$ws.Range("J12").Value = 0.125
$ws.Range("L12").Value = 0.02840909090909091
$ws.Range("S12").Value = 0.05113636363636364
$ws.Range("G13").Value = 0.6428571428571429
$ws.Range("J13").Value = 0.2857142857142857
$ws.Range("S13").Value = 0.07142857142857142
$ws.Range("F15").Value = 0.005617977528089887
$ws.Range("H15").Value = 0.2303370786516854
$ws.Range("I15").Value = 0.06179775280898876
$ws.Range("J15").Value = 0.303370786516854
$ws.Range("K15").Value = 0.101123595505618
$ws.Range("M15").Value = 0.01685393258426966
$ws.Range("O15").Value = 0.06179775280898876
$ws.Range("S15").Value = 0.2191011235955056
$ws.Range("F16").Value = 0.02484472049689441
$ws.Range("H16").Value = 0.1677018633540373
$ws.Range("I16").Value = 0.04347826086956522
$ws.Range("J16").Value = 0.3726708074534161
$ws.Range("K16").Value = 0.1490683229813665
$ws.Range("M16").Value = 0.01863354037267081
$ws.Range("O16").Value = 0.02484472049689441
$ws.Range("S16").Value = 0.1987577639751553
$ws.Range("F17").Value = 0.02064220183486239
$ws.Range("H17").Value = 0.2155963302752294
$ws.Range("I17").Value = 0.06880733944954129
$ws.Range("J17").Value = 0.3967889908256881
$ws.Range("K17").Value = 0.1032110091743119
$ws.Range("M17").Value = 0.01834862385321101
$ws.Range("O17").Value = 0.05045871559633028
$ws.Range("S17").Value = 0.1261467889908257
$ws.Range("F18").Value = 0.01418439716312057
$ws.Range("H18").Value = 0.2056737588652482
$ws.Range("I18").Value = 0.1063829787234043
$ws.Range("J18").Value = 0.4397163120567376
$ws.Range("K18").Value = 0.0851063829787234
$ws.Range("M18").Value = 0.03546099290780142
$ws.Range("O18").Value = 0.04964539007092199
$ws.Range("S18").Value = 0.06382978723404255
$ws.Range("F19").Value = 0.01166666666666667
$ws.Range("H19").Value = 0.2491666666666667
$ws.Range("I19").Value = 0.0525
$ws.Range("J19").Value = 0.3566666666666667
$ws.Range("K19").Value = 0.1116666666666667
$ws.Range("M19").Value = 0.03333333333333333
$ws.Range("N19").Value = 0.0008333333333333334
$ws.Range("O19").Value = 0.07083333333333333
$ws.Range("S19").Value = 0.1133333333333333
